# Update "想去人数" (F column) figures for rows 5-13 (skipping row 9)
# on both the "展览" and "全部类型" worksheets, per the data refresh
# recorded in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 525
    6  = 7075
    7  = 204
    8  = 165
    10 = 432
    11 = 148
    12 = 190
    13 = 612
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
